# 베팅(Bet) 관련 글로벌 상수를 뽑기(Gacha) 관련 글로벌 상수로 변경
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

# Row 8: Bet3Diamonds -> GachaEnergy (value unchanged: 10)
$ws.Range("A8").Value = "GachaEnergy"

# Row 9: Bet3Spins -> Gacha1Event (value 10 -> 1)
$ws.Range("A9").Value = "Gacha1Event"
$ws.Range("B9").Value = 1

# Row 10: Bet3Tickets -> Gacha2Events (value 1 -> 2)
$ws.Range("A10").Value = "Gacha2Events"
$ws.Range("B10").Value = 2

# Row 11: Bet1Event -> Gacha3Events (value 1 -> 10)
$ws.Range("A11").Value = "Gacha3Events"
$ws.Range("B11").Value = 10

# Row 12: Bet2Events -> Gacha1BrokenEnergy (value 2 -> 1)
$ws.Range("A12").Value = "Gacha1BrokenEnergy"
$ws.Range("B12").Value = 1

# Row 13: Bet3Events -> Gacha2BrokenEnergys (value 10 -> 2)
$ws.Range("A13").Value = "Gacha2BrokenEnergys"
$ws.Range("B13").Value = 2
